$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.168.26'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.18%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.852.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.19%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.78%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6847'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.02%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07714'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.86%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3039'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.86%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08160'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.900.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.90%  '

$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7223'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.22%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.192'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.35'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.148.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.20%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007800'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.721'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.00%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '234.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.68%  '

$ws.Range("E21").Value = '  +0.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.099.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.434'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '161.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.55%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.949'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.39%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1428'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.34%  '

$ws.Range("E28").Value = '  -3.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.961'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.55%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.401'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.89%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.520'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.44%  '

$ws.Range("E32").Value = '  -2.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.010'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05186'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.179'
$ws.Range("D35").Style = "Normal"

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7034'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.38%  '

$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.026'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.663'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.50%  '

$ws.Range("E39").Value = '  -4.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.678'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9136'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.100.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.90%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.991'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.09%  '

$ws.Range("E44").Value = '  -4.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.32'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.79%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.757'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.995.33'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.155'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.52%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.901'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.95%  '
